$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The "Prix TTC" totals row (previously row 5: B5/D5/E5, with F5
# trailing along beside it) needs to move down to row 6 to make room
# for a new "Nom de domaine" line item at row 5. F5 itself stays put.
# ------------------------------------------------------------------

# 1) Build the (currently empty) row 6 as the new totals row, copying
#    the number formats from the existing data rows so the same
#    cellXfs entries (font/alignment/numFmt) get reused instead of
#    minting new styles.
$ws.Range("B4").Copy()
$ws.Range("B6").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("D3").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("D6").Value = $ws.Range("D5").Value

$ws.Range("E3").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("E6").Formula = "=(E3+E4+E5)*1.2"

$ws.Rows(6).RowHeight = 15.75

# 2) Now overwrite row 5 (B/D/E) with the new "Nom de domaine" line
#    item - same visual style family as the other data rows (2 & 4).
$ws.Range("A5").Value = "Nom de domaine"

$ws.Range("B4").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("B5").Value = 1

$ws.Range("D3").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("D5").Value = 15

$ws.Range("E3").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$ws.Range("E5").Formula = "=D5*B5"

$excel.CutCopyMode = $false

# 3) Grow "Tableau3" (A2:E5 -> A2:E6) so the new line item is part of
#    the table / AutoFilter range.
$lo1 = $ws.ListObjects.Item(1)
$lo1.Resize($ws.Range("A2:E6"))

# 4) Cosmetic: leave the selection where the author left it on save.
$ws.Range("E7").Select()

Write-Host "edit applied"
